$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "307.67"
Set-TextValue $ws.Range("E2") "0.81%"
Set-TextValue $ws.Range("D3") "41.03"
Set-TextValue $ws.Range("E3") "3.74%"
Set-TextValue $ws.Range("D4") "5.125"
Set-TextValue $ws.Range("E4") "1.78%"
Set-TextValue $ws.Range("E5") "-0.59%"
Set-TextValue $ws.Range("D6") "4.266"
Set-TextValue $ws.Range("E6") "0.36%"
Set-TextValue $ws.Range("D7") "1.617"
Set-TextValue $ws.Range("E7") "1.52%"
Set-TextValue $ws.Range("D9") "0.9010"
Set-TextValue $ws.Range("E9") "2.15%"
Set-TextValue $ws.Range("D10") "0.1094"
Set-TextValue $ws.Range("E10") "12.94%"
Set-TextValue $ws.Range("D11") "0.1772"
Set-TextValue $ws.Range("E11") "3.24%"
Set-TextValue $ws.Range("D12") "0.09169"
Set-TextValue $ws.Range("E12") "3.16%"
Set-TextValue $ws.Range("D13") "0.04201"
Set-TextValue $ws.Range("E13") "-6.39%"
Set-TextValue $ws.Range("E14") "-0.44%"
Set-TextValue $ws.Range("E15") "-1.81%"
Set-TextValue $ws.Range("D16") "0.005804"
Set-TextValue $ws.Range("E16") "-2.42%"
Set-TextValue $ws.Range("E18") "-1.98%"
Set-TextValue $ws.Range("D19") "6.570"
Set-TextValue $ws.Range("E19") "-6.56%"
Set-TextValue $ws.Range("D20") "0.1360"
Set-TextValue $ws.Range("E20") "0.65%"
Set-TextValue $ws.Range("E21") "-12.84%"
Set-TextValue $ws.Range("D22") "0.04072"
Set-TextValue $ws.Range("E22") "-3.06%"
Set-TextValue $ws.Range("D23") "0.001223"
Set-TextValue $ws.Range("E23") "2.37%"
Set-TextValue $ws.Range("D24") "0.004000"
Set-TextValue $ws.Range("E24") "-1.52%"
Set-TextValue $ws.Range("E25") "6.48%"
Set-TextValue $ws.Range("D38") "0.02384"
Set-TextValue $ws.Range("E38") "2.20%"
Set-TextValue $ws.Range("D39") "0.05178"
Set-TextValue $ws.Range("D40") "0.007752"
Set-TextValue $ws.Range("E40") "-2.49%"
Set-TextValue $ws.Range("E41") "-1.59%"
Set-TextValue $ws.Range("D42") "0.006840"
Set-TextValue $ws.Range("E42") "7.74%"
Set-TextValue $ws.Range("D43") "0.001952"
Set-TextValue $ws.Range("E43") "-1.81%"
Set-TextValue $ws.Range("D44") "0.008547"
Set-TextValue $ws.Range("E44") "-1.32%"
Set-TextValue $ws.Range("D45") "0.3074"
Set-TextValue $ws.Range("E45") "1.36%"
Set-TextValue $ws.Range("D46") "0.00006989"
Set-TextValue $ws.Range("E46") "6.88%"
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "-0.07%"
Set-TextValue $ws.Range("D48") "0.02033"
Set-TextValue $ws.Range("E48") "502.29%"
Set-TextValue $ws.Range("D49") "0.004204"
Set-TextValue $ws.Range("E49") "-40.04%"
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "-0.07%"
Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "-0.07%"
